# Updated symbol list on Sat Jan 28 21:43:27 UTC 2023 with GitHub Actions
# Refresh the cryptocurrency price/volume snapshot in columns D (Price) and
# E (Volume(1h)) for the rows whose figures moved since the last run.
# Values are written with a leading apostrophe so Excel stores them as text
# (matching the sheet's existing inline-string cells) rather than
# reinterpreting them as numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.72"
$ws.Range("E2").Value = "'-0.69%"
$ws.Range("D3").Value = "'38.92"
$ws.Range("E3").Value = "'7.39%"
$ws.Range("D4").Value = "'5.099"
$ws.Range("E4").Value = "'0.95%"
$ws.Range("D5").Value = "'0.08082"
$ws.Range("E5").Value = "'-0.43%"
$ws.Range("D6").Value = "'1.924"
$ws.Range("E6").Value = "'-4.25%"
$ws.Range("E7").Value = "'0.87%"
$ws.Range("D8").Value = "'8.041"
$ws.Range("E8").Value = "'2.33%"
$ws.Range("D9").Value = "'0.9274"
$ws.Range("E9").Value = "'0.00%"
$ws.Range("D10").Value = "'0.1459"
$ws.Range("E10").Value = "'-1.74%"
$ws.Range("D11").Value = "'0.1910"
$ws.Range("E11").Value = "'-1.20%"
$ws.Range("D12").Value = "'0.09026"
$ws.Range("E12").Value = "'-1.32%"
$ws.Range("D13").Value = "'0.03508"
$ws.Range("E13").Value = "'-0.67%"
$ws.Range("D14").Value = "'0.09764"
$ws.Range("E14").Value = "'-1.21%"
$ws.Range("D15").Value = "'0.001401"
$ws.Range("E15").Value = "'-0.54%"
$ws.Range("D16").Value = "'0.005867"
$ws.Range("E16").Value = "'-3.88%"
$ws.Range("D17").Value = "'3.784"
$ws.Range("E17").Value = "'-1.53%"
$ws.Range("D18").Value = "'3.409"
$ws.Range("E18").Value = "'-0.79%"
$ws.Range("D19").Value = "'0.3462"
$ws.Range("E19").Value = "'0.13%"
$ws.Range("D20").Value = "'0.1327"
$ws.Range("E20").Value = "'2.62%"
$ws.Range("D21").Value = "'4.693"
$ws.Range("E21").Value = "'-2.80%"
$ws.Range("D23").Value = "'0.04377"
$ws.Range("E23").Value = "'-0.21%"
$ws.Range("D24").Value = "'0.001234"
$ws.Range("E24").Value = "'-0.15%"
$ws.Range("D25").Value = "'0.004274"
$ws.Range("E25").Value = "'2.06%"
$ws.Range("D26").Value = "'0.0001303"
$ws.Range("E26").Value = "'0.18%"
$ws.Range("D39").Value = "'0.02031"
$ws.Range("E39").Value = "'-1.25%"
$ws.Range("E40").Value = "'-1.62%"
$ws.Range("D41").Value = "'0.007523"
$ws.Range("E41").Value = "'0.69%"
$ws.Range("D42").Value = "'0.009910"
$ws.Range("E42").Value = "'-0.86%"
$ws.Range("D43").Value = "'0.1342"
$ws.Range("E43").Value = "'-2.26%"
$ws.Range("D44").Value = "'0.002105"
$ws.Range("E44").Value = "'-0.77%"
$ws.Range("D45").Value = "'0.009929"
$ws.Range("E45").Value = "'0.42%"
$ws.Range("D46").Value = "'0.00006202"
$ws.Range("E46").Value = "'-1.56%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.09%"
$ws.Range("D48").Value = "'0.002875"
$ws.Range("D49").Value = "'0.001803"
$ws.Range("E49").Value = "'12.64%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.09%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.09%"
